# Sprint 39 test-case-summary update:
#  - Day 5 ("Total testcase Written") execution count bumped 912 -> 920
#  - Day 6 block (rows 33-35) gets its first real numbers: Written/Execution/Review
#  - Day 7 block (rows 39-41) gets its first real numbers: Written/Execution/Review
#  - The view scrolls down and the selection moves to C41 (last cell touched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 5 (rows 26-29): bump the "Total testcase Written" count
$ws.Range("C27").Value = 920

# Day 6 (rows 32-35): fill in Written / Execution / Review totals
$ws.Range("C33").Value = 948
$ws.Range("C34").Value = 1172
$ws.Range("C35").Value = 636

# Day 7 (rows 38-41): fill in Written / Execution / Review totals
$ws.Range("C39").Value = 983
$ws.Range("C40").Value = 1222
$ws.Range("C41").Value = 686

# Move the visible window / selection down to where the new data was entered
$ws.Range("C41").Select()
